# Apply the edits described by the commit:
#  1. Rename header text "Prob. of Failure" -> "Prob_of_Failure"
#  2. Change the font color used by the numeric "Prob. of Failure" column
#     (rows 2-19, column C) from the automatic theme color to explicit black
#  3. Normalize the custom row heights to 19.5 for all data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the header cell text in C1
$ws.Range("C1").Value = "Prob_of_Failure"

# 2. Make the font color of the data cells (C2:C19) explicit black
$ws.Range("C2:C19").Font.Color = 0

# 3. Set a uniform row height of 19.5 for every used row (1-19)
$ws.Range("A1:C19").RowHeight = 19.5
